$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datenerfassung")

# --- Fill in the previously-empty time-record rows 59-71 -------------------
# Columns: A = Datum (date serial), B = Dauer (fraction of a day),
#          C = Primäre Tätigkeit (picklist text), D = Anmerkung / LV-Einheit
# Values are written top-to-bottom, left-to-right so any brand-new picklist
# text lands in the shared-string table in the same order it was typed.

$ws.Cells.Item(59, 1).Value = 44326
$ws.Cells.Item(59, 2).Value = 2 / 24
$ws.Cells.Item(59, 3).Value = "LV-Einheit"

$ws.Cells.Item(60, 1).Value = 44327
$ws.Cells.Item(60, 2).Value = 1 / 24
$ws.Cells.Item(60, 3).Value = "Koordination und Projektmanagement"
$ws.Cells.Item(60, 4).Value = "Teammeeting"

$ws.Cells.Item(61, 1).Value = 44328
$ws.Cells.Item(61, 2).Value = 3 / 24
$ws.Cells.Item(61, 3).Value = "Implementierung"
$ws.Cells.Item(61, 4).Value = "Issues Nr. 71, 83, 86"

$ws.Cells.Item(62, 1).Value = 44328
$ws.Cells.Item(62, 2).Value = 0.5 / 24
$ws.Cells.Item(62, 3).Value = "Konfiguration und Deployment"
$ws.Cells.Item(62, 4).Value = "Docker Lösung testen mit Würfel"

$ws.Cells.Item(63, 1).Value = 44329
$ws.Cells.Item(63, 2).Value = 2 / 24
$ws.Cells.Item(63, 3).Value = "Koordination und Projektmanagement"
$ws.Cells.Item(63, 4).Value = "Teammeeting"

$ws.Cells.Item(64, 1).Value = 44329
$ws.Cells.Item(64, 2).Value = 3 / 24
$ws.Cells.Item(64, 3).Value = "Koordination und Projektmanagement"
$ws.Cells.Item(64, 4).Value = "Teammeeting"

$ws.Cells.Item(65, 1).Value = 44333
$ws.Cells.Item(65, 2).Value = 1 / 24
$ws.Cells.Item(65, 3).Value = "Systemtest (fremdes System)"
$ws.Cells.Item(65, 4).Value = "Installation fremdes System"

$ws.Cells.Item(66, 1).Value = 44333
$ws.Cells.Item(66, 2).Value = 1 / 24
$ws.Cells.Item(66, 3).Value = "Systemtest (fremdes System)"

$ws.Cells.Item(67, 1).Value = 44335
$ws.Cells.Item(67, 2).Value = 2 / 24
$ws.Cells.Item(67, 3).Value = "Systemtest (fremdes System)"

$ws.Cells.Item(68, 1).Value = 44335
$ws.Cells.Item(68, 2).Value = 2 / 24
$ws.Cells.Item(68, 3).Value = "Systemtest (fremdes System)"
$ws.Cells.Item(68, 4).Value = "Installation fremdes System"

$ws.Cells.Item(69, 1).Value = 44336
$ws.Cells.Item(69, 2).Value = 3 / 24
$ws.Cells.Item(69, 3).Value = "Systemtest (fremdes System)"

$ws.Cells.Item(70, 1).Value = 44343
$ws.Cells.Item(70, 2).Value = 1 / 24
$ws.Cells.Item(70, 3).Value = "Koordination und Projektmanagement"
$ws.Cells.Item(70, 4).Value = "Teammeeting"

$ws.Cells.Item(71, 1).Value = 44345
$ws.Cells.Item(71, 2).Value = 1.5 / 24
$ws.Cells.Item(71, 3).Value = "Implementierung"
$ws.Cells.Item(71, 4).Value = "Issue Nr. 129"

# --- Make room for one more row before the trailing blank row --------------
# Row 1016 used to be the final (empty, differently-styled) sentinel row;
# insert a fresh blank data row above it (copying the data rows' style) so
# the table can keep growing and the old sentinel becomes row 1017.
$ws.Rows.Item(1016).Insert()

# --- Restore the scroll / selection state shown in the saved workbook ------
$ws.Activate()
$ws.Range("C2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B72").Select()
